# Set "想去人数" (column F) to 0 for all data rows on the
# relevant worksheets, leaving the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "演出", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 6).Value = 0
    }
}
